# Auto-generated edit script applying the recorded cell-value changes
# from the commit diff. All target cells are plain numeric literals
# (no formulas in the source workbook), so we just set/clear .Value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 385.31818
$ws.Range("I28").Value = 426.84616
$ws.Range("K28").Value = 426.84616
$ws.Range("M28").Value = 58.15384

$ws.Range("H62").Value = 3277.4443
$ws.Range("I62").Value = 3062.125
$ws.Range("K62").Value = 3062.125
$ws.Range("M62").Value = -2438.125

$ws.Range("H65").Value = 3277.4443
$ws.Range("I65").Value = 3062.125
$ws.Range("K65").Value = 15310.625
$ws.Range("M65").Value = -12190.625

$ws.Range("H94").Value = 866.5625
$ws.Range("I94").Value = 866.5625
$ws.Range("K94").Value = 866.5625
$ws.Range("M94").Value = -415.5625

$ws.Range("H112").Value = 4176.2085
$ws.Range("J112").Value = 4516.45
$ws.Range("L112").Value = 13549.35
$ws.Range("N112").Value = -15765.35

$ws.Range("H127").Value = 1034.3334
$ws.Range("I127").Value = 1034.3334
$ws.Range("K127").Value = 3103.0002
$ws.Range("M127").Value = 1856.9998

$ws.Range("H141").Value = 3758.8
$ws.Range("I141").Value = 3060.1538
$ws.Range("K141").Value = 9180.4614
$ws.Range("M141").Value = -4000.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 162706.08
$ws.Range("I32").Value = 173303.69
$ws.Range("J32").Value = 9040.75
$ws.Range("K32").Value = 173303.69
$ws.Range("L32").Value = 9040.75
$ws.Range("M32").Value = -173016.69
$ws.Range("N32").Value = -9614.75

$ws.Range("H45").Value = 252921.5
$ws.Range("I45").Value = 252921.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 252921.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -252544.5
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 967104.25
$ws.Range("I61").Value = 2682.48
$ws.Range("J61").Value = 5350839.5
$ws.Range("K61").Value = 2682.48
$ws.Range("L61").Value = 5350839.5
$ws.Range("M61").Value = -2470.48
$ws.Range("N61").Value = -5351263.5

$ws.Range("H102").Value = 13006.889
$ws.Range("I102").Value = 18495
$ws.Range("J102").Value = 2030.6666
$ws.Range("K102").Value = 18495
$ws.Range("L102").Value = 2030.6666
$ws.Range("M102").Value = -16873
$ws.Range("N102").Value = -5274.6666

$ws.Range("H110").Value = 2097.1177
$ws.Range("I110").Value = 1867.9286
$ws.Range("K110").Value = 1867.9286
$ws.Range("M110").Value = 177.0714

$ws.Range("H127").Value = 129990
$ws.Range("J127").Value = 129990
$ws.Range("L127").Value = 129990
$ws.Range("N127").Value = -139910

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H136").Value = 967104.25
$ws.Range("I136").Value = 2682.48
$ws.Range("J136").Value = 5350839.5
$ws.Range("K136").Value = 8047.440000000001
$ws.Range("L136").Value = 16052518.5
$ws.Range("M136").Value = -5497.440000000001
$ws.Range("N136").Value = -16057618.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 3196.5
$ws.Range("I29").Value = 3196.5
$ws.Range("K29").Value = 3196.5
$ws.Range("M29").Value = -2907.5

$ws.Range("H86").Value = 2993.4666
$ws.Range("I86").Value = 1090.5
$ws.Range("K86").Value = 1090.5
$ws.Range("M86").Value = 32.5

$ws.Range("H89").Value = 2993.4666
$ws.Range("I89").Value = 1090.5
$ws.Range("K89").Value = 5452.5
$ws.Range("M89").Value = 163.5

$ws.Range("H94").Value = 15745.923
$ws.Range("I94").Value = 18076.555
$ws.Range("K94").Value = 18076.555
$ws.Range("M94").Value = -17625.555

$ws.Range("H105").Value = 9659.532999999999
$ws.Range("I105").Value = 12821.777
$ws.Range("J105").Value = 4916.1665
$ws.Range("K105").Value = 12821.777
$ws.Range("L105").Value = 4916.1665
$ws.Range("M105").Value = -11074.777
$ws.Range("N105").Value = -8410.166499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 189.5
$ws.Range("I7").Value = 214.4
$ws.Range("J7").Value = 65
$ws.Range("K7").Value = 214.4
$ws.Range("L7").Value = 65
$ws.Range("M7").Value = -101.4
$ws.Range("N7").Value = -291

$ws.Range("H31").Value = 2760.1128
$ws.Range("I31").Value = 2541.0908
$ws.Range("K31").Value = 2541.0908
$ws.Range("M31").Value = -2246.0908

$ws.Range("H34").Value = 2760.1128
$ws.Range("I34").Value = 2541.0908
$ws.Range("K34").Value = 2541.0908
$ws.Range("M34").Value = -2339.0908

$ws.Range("H58").Value = 2847.3333
$ws.Range("I58").Value = 2482.35
$ws.Range("J58").Value = 3577.3
$ws.Range("K58").Value = 2482.35
$ws.Range("L58").Value = 3577.3
$ws.Range("M58").Value = -2279.35
$ws.Range("N58").Value = -3983.3

$ws.Range("H132").Value = 30109.594
$ws.Range("I132").Value = 50516.855
$ws.Range("J132").Value = 3325.0625
$ws.Range("K132").Value = 151550.565
$ws.Range("L132").Value = 9975.1875
$ws.Range("M132").Value = -149020.565
$ws.Range("N132").Value = -15035.1875

$ws.Range("H134").Value = 1334.5349
$ws.Range("I134").Value = 1343.7241
$ws.Range("K134").Value = 4031.1723
$ws.Range("M134").Value = -1496.1723

$ws.Range("H136").Value = 2847.3333
$ws.Range("I136").Value = 2482.35
$ws.Range("J136").Value = 3577.3
$ws.Range("K136").Value = 7447.049999999999
$ws.Range("L136").Value = 10731.9
$ws.Range("M136").Value = -4897.049999999999
$ws.Range("N136").Value = -15831.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4129362
$ws.Range("I4").Value = 3583605.2
$ws.Range("K4").Value = 10750815.6
$ws.Range("M4").Value = -10750703.6

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H86").Value = 1673.6
$ws.Range("J86").Value = 1727.5714
$ws.Range("L86").Value = 5182.7142
$ws.Range("N86").Value = -7554.7142

$ws.Range("H89").Value = 1673.6
$ws.Range("J89").Value = 1727.5714
$ws.Range("L89").Value = 15548.1426
$ws.Range("N89").Value = -27404.1426

$ws.Range("H104").Value = 9478.125
$ws.Range("J104").Value = 13499
$ws.Range("L104").Value = 40497
$ws.Range("N104").Value = -45739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4352.037
$ws.Range("I97").Value = 652.9
$ws.Range("K97").Value = 652.9
$ws.Range("M97").Value = -156.9

$ws.Range("H107").Value = 43944.695
$ws.Range("I107").Value = 71595.57000000001
$ws.Range("K107").Value = 71595.57000000001
$ws.Range("M107").Value = -69675.57000000001

$ws.Range("H132").Value = 515288.6
$ws.Range("I132").Value = 2011.375
$ws.Range("K132").Value = 6034.125
$ws.Range("M132").Value = -3504.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2226.8572
$ws.Range("I16").Value = 2090.5
$ws.Range("K16").Value = 2090.5
$ws.Range("M16").Value = -1920.5

$ws.Range("H36").Value = 68000
$ws.Range("J36").Value = 68000
$ws.Range("L36").Value = 68000
$ws.Range("N36").Value = -69124

$ws.Range("H40").Value = 3973.8125
$ws.Range("I40").Value = 2359
$ws.Range("J40").Value = 5588.625
$ws.Range("K40").Value = 2359
$ws.Range("L40").Value = 5588.625
$ws.Range("M40").Value = -2223
$ws.Range("N40").Value = -5860.625

$ws.Range("H93").Value = 2112.4614
$ws.Range("I93").Value = 1593.2354
$ws.Range("J93").Value = 3093.2222
$ws.Range("K93").Value = 1593.2354
$ws.Range("L93").Value = 3093.2222
$ws.Range("M93").Value = -345.2354
$ws.Range("N93").Value = -5589.2222

$ws.Range("H132").Value = 3033.3333
$ws.Range("I132").Value = 2862.3215
$ws.Range("K132").Value = 8586.9645
$ws.Range("M132").Value = -6056.9645

$ws.Range("H136").Value = 3521.6458
$ws.Range("I136").Value = 2467.2
$ws.Range("J136").Value = 4667.7827
$ws.Range("K136").Value = 7401.599999999999
$ws.Range("L136").Value = 14003.3481
$ws.Range("M136").Value = -4851.599999999999
$ws.Range("N136").Value = -19103.3481

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 19475
$ws.Range("I74").Value = 10000
$ws.Range("K74").Value = 10000
$ws.Range("M74").Value = -9064

$ws.Range("H77").Value = 19475
$ws.Range("I77").Value = 10000
$ws.Range("K77").Value = 30000
$ws.Range("M77").Value = -25320

$ws.Range("H96").Value = 21334.8
$ws.Range("I96").Value = 1668.5
$ws.Range("J96").Value = 100000
$ws.Range("K96").Value = 1668.5
$ws.Range("L96").Value = 100000
$ws.Range("M96").Value = -295.5
$ws.Range("N96").Value = -102746

$ws.Range("H113").Value = 933
$ws.Range("J113").Value = 933
$ws.Range("L113").Value = 2799
$ws.Range("N113").Value = -7139

$ws.Range("H122").Value = 1964.4474
$ws.Range("I122").Value = 1876.6
$ws.Range("K122").Value = 5629.799999999999
$ws.Range("M122").Value = -3179.799999999999

$ws.Range("H124").Value = 87750
$ws.Range("J124").Value = 87750
$ws.Range("L124").Value = 87750
$ws.Range("N124").Value = -97570

